$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.014.80"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").Value = "3.065.12"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'516.61"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "'141.48"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").Value = "'7.30"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "3.592.23"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").Value = "'26.42"
$ws.Range("E14").Value = "  +4.70%  "
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "58.052.05"
$ws.Range("D17").Value = "3.070.89"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").Value = "'6.12"
$ws.Range("E18").Value = "  +2.93%  "
$ws.Range("D19").Value = "'12.83"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("D20").Value = "'8.19"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").Value = "'331.24"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'0.501"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'65.37"
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "0.0₃0908"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").Value = "'6.47"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "'7.23"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  +3.66%  "
$ws.Range("D32").Value = "'20.63"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").Value = "'155.17"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").Value = "'27.69"
$ws.Range("E34").Value = "  +2.55%  "
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("D38").Value = "'0.0678"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").Value = "3.108.64"
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("D40").Value = "'3.91"
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("D41").Value = "'36.81"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "2.299.15"
$ws.Range("E44").Value = "  +4.50%  "
$ws.Range("E45").Value = "  +4.15%  "
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("D47").Value = "'20.88"
$ws.Range("E47").Value = "  +5.64%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'5.94"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'0.938"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("E50").Value = "  +10.02%  "
$ws.Range("D51").Value = "'253.93"
$ws.Range("E51").Value = "  +9.86%  "
